# Refresh the cryptocurrency ranking table (coin name/link/price/1h volume)
# for the rows whose data moved since the last GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "42.266.03"
$ws.Cells.Item(2, 5).Value = "  -0.97%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "2.274.79"
$ws.Cells.Item(3, 5).Value = "  -1.31%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "'299.95"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.35%  "

# Row 6: Solana
$ws.Cells.Item(6, 4).Value = "'96.43"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.80%  "

# Row 7: XRP
$ws.Cells.Item(7, 4).Value = "'0.497"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -1.71%  "

# Row 8: USDC
$ws.Cells.Item(8, 5).Value = "  +0.03%  "

# Row 9: Cardano
$ws.Cells.Item(9, 4).Value = "'0.494"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.77%  "

# Row 10: Avalanche
$ws.Cells.Item(10, 4).Value = "'33.36"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -3.11%  "

# Row 11: Dogecoin
$ws.Cells.Item(11, 5).Value = "  -0.15%  "

# Row 12: OKB
$ws.Cells.Item(12, 4).Value = "'48.07"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -7.15%  "

# Row 13: TRON
$ws.Cells.Item(13, 5).Value = "  +0.20%  "

# Row 14: Polkadot
$ws.Cells.Item(14, 4).Value = "'6.68"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.11%  "

# Row 15: Chainlink
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "'15.60"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.56%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Cells.Item(16, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(16, 4).Value = "2.623.42"
$ws.Cells.Item(16, 5).Value = "  -1.64%  "

# Row 17: WrappedEther
$ws.Cells.Item(17, 4).Value = "2.266.43"
$ws.Cells.Item(17, 5).Value = "  -1.76%  "

# Row 18: Polygon
$ws.Cells.Item(18, 4).Value = "'0.786"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -4.44%  "

# Row 19: WrappedBTC
$ws.Cells.Item(19, 4).Value = "42.166.14"
$ws.Cells.Item(19, 5).Value = "  -1.07%  "

# Row 20: InternetComputer(DFINITY)
$ws.Cells.Item(20, 4).Value = "'11.71"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.40%  "

# Row 21: ShibaInu
$ws.Cells.Item(21, 4).Value = "0.0₃0892"
$ws.Cells.Item(21, 5).Value = "  -1.06%  "

# Row 22: Uniswap
$ws.Cells.Item(22, 4).Value = "'6.00"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.76%  "

# Row 23: Litecoin
$ws.Cells.Item(23, 4).Value = "'66.52"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.44%  "

# Row 24: BitcoinCash
$ws.Cells.Item(24, 4).Value = "'235.30"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.21%  "

# Row 25: ImmutableX
$ws.Cells.Item(25, 4).Value = "'1.99"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.44%  "

# Row 26: PancakeSwap
$ws.Cells.Item(26, 2).Value = "PancakeSwap"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(26, 4).Value = "'2.46"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.22%  "

# Row 27: Dai
$ws.Cells.Item(27, 2).Value = "Dai"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "

# Row 28: EthereumClassic
$ws.Cells.Item(28, 4).Value = "'24.03"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -4.65%  "

# Row 29: Monero
$ws.Cells.Item(29, 2).Value = "Monero"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(29, 4).Value = "'168.88"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.63%  "

# Row 30: Toncoin
$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).Value = "'2.07"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.42%  "

# Row 31: Cosmos
$ws.Cells.Item(31, 4).Value = "'9.20"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.10%  "

# Row 32: InjectiveProtocol
$ws.Cells.Item(32, 4).Value = "'33.75"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.54%  "

# Row 33: FirstDigitalUSD
$ws.Cells.Item(33, 4).Value = "'0.999"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.11%  "

# Row 34: Filecoin
$ws.Cells.Item(34, 4).Value = "'4.91"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.64%  "

# Row 35: RenderToken
$ws.Cells.Item(35, 4).Value = "'4.51"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -2.49%  "

# Row 36: Celestia
$ws.Cells.Item(36, 2).Value = "Celestia"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(36, 4).Value = "'16.68"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.52%  "

# Row 37: WEMIXToken
$ws.Cells.Item(37, 2).Value = "WEMIXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(37, 4).Value = "'2.33"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.69%  "

# Row 38: Hedera
$ws.Cells.Item(38, 4).Value = "'0.0687"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -4.15%  "

# Row 39: LidoDAOToken
$ws.Cells.Item(39, 4).Value = "'2.79"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.57%  "

# Row 40: Kaspa
$ws.Cells.Item(40, 4).Value = "'0.0989"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.49%  "

# Row 41: Stellar
$ws.Cells.Item(41, 5).Value = "  -2.40%  "

# Row 42: ARBITRUM
$ws.Cells.Item(42, 5).Value = "  -4.66%  "

# Row 43: ApeXProtocol
$ws.Cells.Item(43, 4).Value = "'2.42"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.13%  "

# Row 44: Maker
$ws.Cells.Item(44, 4).Value = "1.971.13"
$ws.Cells.Item(44, 5).Value = "  -1.05%  "

# Row 45: VeChain
$ws.Cells.Item(45, 4).Value = "'0.0278"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.11%  "

# Row 46: EnergySwap
$ws.Cells.Item(46, 4).Value = "'17.43"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -7.09%  "

# Row 47: FraxShare
$ws.Cells.Item(47, 4).Value = "'9.56"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -6.57%  "

# Row 48: NEARProtocol
$ws.Cells.Item(48, 4).Value = "'2.79"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.69%  "

# Row 49: RocketPoolETH
$ws.Cells.Item(49, 4).Value = "2.496.88"
$ws.Cells.Item(49, 5).Value = "  -1.59%  "

# Row 50: MultiversX
$ws.Cells.Item(50, 4).Value = "'52.34"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -5.29%  "

# Row 51: Stacks
$ws.Cells.Item(51, 2).Value = "Stacks"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(51, 4).Value = "'1.48"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.66%  "
